# Apply scheduled profit-column updates to Aegis_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 125017
$ws.Range("I11").Value = 125017
$ws.Range("K11").Value = 125017
$ws.Range("M11").Value = -124877

$ws.Range("H51").Value = 8104.7617
$ws.Range("J51").Value = 3242.7856
$ws.Range("L51").Value = 3242.7856
$ws.Range("N51").Value = -4210.7856

$ws.Range("H55").Value = 204.77777
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 234.71428
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 234.71428
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -662.71428

$ws.Range("H86").Value = 2166
$ws.Range("I86").Value = 2101
$ws.Range("J86").Value = 2382.6667
$ws.Range("K86").Value = 2101
$ws.Range("L86").Value = 2382.6667
$ws.Range("M86").Value = -978
$ws.Range("N86").Value = -4628.6667

$ws.Range("H89").Value = 2166
$ws.Range("I89").Value = 2101
$ws.Range("J89").Value = 2382.6667
$ws.Range("K89").Value = 10505
$ws.Range("L89").Value = 11913.3335
$ws.Range("M89").Value = -4889
$ws.Range("N89").Value = -23145.3335

$ws.Range("H92").Value = 514.8333
$ws.Range("I92").Value = 455.2857
$ws.Range("K92").Value = 455.2857
$ws.Range("M92").Value = 792.7143

$ws.Range("H112").Value = 1045.8148
$ws.Range("J112").Value = 1066.8077
$ws.Range("L112").Value = 3200.4231
$ws.Range("N112").Value = -5416.4231

$ws.Range("H121").Value = 963.5714
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 883.8461
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 2651.5383
$ws.Range("M121").Value = -4253
$ws.Range("N121").Value = -6145.5383

$ws.Range("H129").Value = 3061.818
$ws.Range("I129").Value = 7168.2
$ws.Range("J129").Value = 937.8276
$ws.Range("K129").Value = 21504.6
$ws.Range("L129").Value = 2813.4828
$ws.Range("M129").Value = -16504.6
$ws.Range("N129").Value = -12813.4828

$ws.Range("H138").Value = 2349.7324
$ws.Range("I138").Value = 1804.6875
$ws.Range("J138").Value = 2508.291
$ws.Range("K138").Value = 5414.0625
$ws.Range("L138").Value = 7524.873000000001
$ws.Range("M138").Value = -274.0625
$ws.Range("N138").Value = -17804.873

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3562.68
$ws.Range("I32").Value = 3292.4536
$ws.Range("J32").Value = 12300
$ws.Range("K32").Value = 3292.4536
$ws.Range("L32").Value = 12300
$ws.Range("M32").Value = -3005.4536
$ws.Range("N32").Value = -12874

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2258.4688
$ws.Range("I20").Value = 2485.9583
$ws.Range("J20").Value = 1576
$ws.Range("K20").Value = 2485.9583
$ws.Range("L20").Value = 1576
$ws.Range("M20").Value = -2238.9583
$ws.Range("N20").Value = -2070

$ws.Range("H64").Value = 1180.0834
$ws.Range("J64").Value = 1361.8889
$ws.Range("L64").Value = 1361.8889
$ws.Range("N64").Value = -1811.8889

$ws.Range("H67").Value = 1180.0834
$ws.Range("J67").Value = 1361.8889
$ws.Range("L67").Value = 1361.8889
$ws.Range("N67").Value = -2921.8889

$ws.Range("H107").Value = 90909980
$ws.Range("I107").Value = 125000860
$ws.Range("J107").Value = 967
$ws.Range("K107").Value = 125000860
$ws.Range("L107").Value = 967
$ws.Range("M107").Value = -124998940
$ws.Range("N107").Value = -4807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38992.13
$ws.Range("I31").Value = 733.1818
$ws.Range("J31").Value = 54022.43
$ws.Range("K31").Value = 733.1818
$ws.Range("L31").Value = 54022.43
$ws.Range("M31").Value = -438.1818
$ws.Range("N31").Value = -54612.43

$ws.Range("H34").Value = 38992.13
$ws.Range("I34").Value = 733.1818
$ws.Range("J34").Value = 54022.43
$ws.Range("K34").Value = 733.1818
$ws.Range("L34").Value = 54022.43
$ws.Range("M34").Value = -531.1818
$ws.Range("N34").Value = -54426.43

$ws.Range("H80").Value = 11383
$ws.Range("J80").Value = 11383
$ws.Range("L80").Value = 11383
$ws.Range("N80").Value = -13629

$ws.Range("H83").Value = 11383
$ws.Range("J83").Value = 11383
$ws.Range("L83").Value = 34149
$ws.Range("N83").Value = -45381

$ws.Range("H132").Value = 4720.625
$ws.Range("I132").Value = 4876.1816
$ws.Range("J132").Value = 4378.4
$ws.Range("K132").Value = 14628.5448
$ws.Range("L132").Value = 13135.2
$ws.Range("M132").Value = -12098.5448
$ws.Range("N132").Value = -18195.2

$ws.Range("H134").Value = 2679.1428
$ws.Range("I134").Value = 2976.8572
$ws.Range("J134").Value = 2381.4285
$ws.Range("K134").Value = 8930.571599999999
$ws.Range("L134").Value = 7144.2855
$ws.Range("M134").Value = -6395.571599999999
$ws.Range("N134").Value = -12214.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1080
$ws.Range("J34").Value = 1300
$ws.Range("L34").Value = 3900
$ws.Range("N34").Value = -4068

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H122").Value = 603.04346
$ws.Range("I122").Value = 603.4
$ws.Range("J122").Value = 602.94446
$ws.Range("K122").Value = 5430.599999999999
$ws.Range("L122").Value = 5426.50014
$ws.Range("M122").Value = -2980.599999999999
$ws.Range("N122").Value = -10326.50014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142863260
$ws.Range("I80").Value = 250009220
$ws.Range("K80").Value = 250009220
$ws.Range("M80").Value = -250008222

$ws.Range("H83").Value = 142863260
$ws.Range("I83").Value = 250009220
$ws.Range("K83").Value = 1250046100
$ws.Range("M83").Value = -1250041108

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 482538.1
$ws.Range("I46").Value = 298.75
$ws.Range("K46").Value = 298.75
$ws.Range("M46").Value = -110.75

$ws.Range("H121").Value = 29470
$ws.Range("J121").Value = 29470
$ws.Range("L121").Value = 29470
$ws.Range("N121").Value = -32964

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1319.8657
$ws.Range("I136").Value = 431.11765
$ws.Range("J136").Value = 2235.5454
$ws.Range("K136").Value = 1293.35295
$ws.Range("L136").Value = 6706.6362
$ws.Range("M136").Value = 1256.64705
$ws.Range("N136").Value = -11806.6362
